# Applies the "add io-ts for api response validation" update to the
# fitness-app tuntikirjanpito (hours log) workbook:
#  - B11: 3 -> 5 hours
#  - B12: 5 -> 6 hours
#  - new row 13: 28.3.2019 / 7 hours / long description of the day's work
#  - sheet scrolled so topLeftCell is A3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hour totals on the two existing rows that changed.
$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 6

# Add the new log entry as row 13, copying the formatting (styles,
# wrap text, alignment, number format) from row 12 so the new row
# matches the look of the rest of the table.
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A13").Value = "28.3.2019"
$ws.Range("B13").Value = 7
$ws.Range("C13").Value = "json-server backendiin, frontin service aloitus, api vastausten TS-validoinnin selvittelyä, validointi io-ts -kirjastolla, sen opettelua, karmea määrä TS -tyyppien refaktorointia käyttämään io-ts tyyppejä"

# Row grew tall because of the wrapped description text.
$ws.Rows.Item(13).RowHeight = 73.2

# Scroll position changed in the saved view.
$ws.Application.ActiveWindow.ScrollRow = 3
